$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the Attribute (A) and Type (B) columns, rows 2-21,
# reflecting the reordered/combined sensor data attribute list.
$data = @(
    @("concept:name", "str"),
    @("operation_end_time", "datetime"),
    @("lifecycle:transition", "str"),
    @("current_task", "str"),
    @("complete_service_time", "str"),
    @("time:timestamp", "datetime"),
    @("process_model_id", "str"),
    @("planned_operation_time", "str"),
    @("SubProcessID", "str"),
    @("case:concept:name", "str"),
    @("case", "str"),
    @("identifier:id", "str"),
    @("requested_service_url", "str"),
    @("parameters", "dict"),
    @("unsatisfied_condition_description", "str"),
    @("org:resource", "str"),
    @("response_status_code", "float"),
    @("event_id", "str"),
    @("human_workstation_green_button_pressed", "float"),
    @("lifecycle:state", "str")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
